$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1. Title
# ---------------------------------------------------------------
$d.Content.Find.Execute(
    "Nanobots: The Future of Medicine", $true, $false, $false, $false, $false,
    $true, 1, $false, "The Heartbeat of Healing: A Journey into Medicine", 2
) | Out-Null

# ---------------------------------------------------------------
# 2. Author line: "Dr. Katherine Abernathy" -> "Amelia Harrison"
# ---------------------------------------------------------------
$d.Content.Find.Execute(
    "Dr. Katherine Abernathy", $true, $false, $false, $false, $false,
    $true, 1, $false, "Amelia Harrison", 2
) | Out-Null

# ---------------------------------------------------------------
# 3. Email line: "katherine.abernathy@healthcare.edu" -> "ameliaharr1021@gmail.com"
# ---------------------------------------------------------------
$d.Content.Find.Execute(
    "katherine", $true, $false, $false, $false, $false,
    $true, 1, $false, "ameliaharr1021@gmail", 2
) | Out-Null

$d.Content.Find.Execute(
    "abernathy@healthcare.edu", $true, $false, $false, $false, $false,
    $true, 1, $false, "com", 2
) | Out-Null

# ---------------------------------------------------------------
# 4. Big intro paragraph -- replace the whole body with new text,
#    using Chr(11) (vertical tab) for the in-paragraph line breaks
#    (<w:br/>) that separate the "sections" of the essay.
# ---------------------------------------------------------------
$break = [char]11
$bigText = (
    "As we embark on this chapter of learning and exploration, we will delve into the realm of medicine, an extraordinary field dedicated to preserving and restoring human health." +
    " From the intricate harmony of our physiology to the boundless possibilities of modern therapies, we will discover the beauty and complexities of medicine." +
    $break + $break +
    "The human body, a marvel of biological engineering, serves as the canvas upon which medicine works its magic." +
    " We will journey through the interconnected systems that govern our existence, from the microscopic world of cells to the intricate network of organs." +
    " We will unravel the secrets of human physiology, discovering the delicate balance that maintains our health and the ways in which medicine can intervene when this balance is disrupted." +
    $break + $break +
    "At the heart of medicine lies the patient, an individual with unique stories, hopes, and fears." +
    " As we step into the world of healing, we will explore the art of patient care, emphasizing empathy, compassion, and respect for human dignity." +
    " We will learn how medicine goes beyond treating symptoms; it also involves nurturing the human spirit, fostering a bond of trust that empowers patients to actively participate in their healing journey." +
    $break + $break +
    "Introduction Continued:" +
    $break + $break +
    "The practice of medicine draws upon a vast reservoir of knowledge, including scientific research, clinical experience, and cultural traditions." +
    " We will examine how evidence-based medicine guides medical practice, ensuring that treatments are safe, effective, and tailored to individual patient needs." +
    " We will also explore alternative and complementary therapies, recognizing the diversity of approaches to healing." +
    $break + $break +
    "As medicine evolves, it faces numerous challenges, from emerging diseases and antibiotic resistance to the complexities of healthcare systems and the rising cost of medical care." +
    " We will delve into these issues, seeking a deeper understanding of the intricate factors that shape modern medicine." +
    " Through critical thinking and open-minded discussions, we will explore potential solutions and envision a future where medical advancements benefit all members of society." +
    $break + $break +
    "Introduction Concluded:" +
    $break + $break +
    "Medicine is a noble profession, one that intertwines science, art, and unwavering dedication to patient care." +
    " As we traverse the landscape of healing, may we cultivate a deep appreciation for the complex tapestry of human life, the indomitable power of the human spirit, and the profound responsibility we carry as future guardians of health."
)

$bigPara = $d.Paragraphs.Item(5)
$bigRange = $d.Range($bigPara.Range.Start, $bigPara.Range.End)
$bigRange.Text = $bigText

# ---------------------------------------------------------------
# 5. Summary paragraph
# ---------------------------------------------------------------
$sumText = (
    "Our exploration of medicine has unveiled the profound impact it wields on human lives, delving into the remarkable complexity of the human body, the challenges of modern medicine, and the essential role of empathy and respect in patient care." +
    " We have unraveled the delicate interplay between science, tradition, and cultural factors that shape medical practices." +
    " Ultimately, we have gained an appreciation for the profound responsibility that rests upon those dedicated to preserving and restoring human health."
)

$sumPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$sumRange = $d.Range($sumPara.Range.Start, $sumPara.Range.End)
$sumRange.Text = $sumText

# ---------------------------------------------------------------
# 6. Append a trailing empty paragraph at the very end of the body.
# ---------------------------------------------------------------
$endPos = $d.Content.End - 1
$endRange = $d.Range($endPos, $endPos)
$endRange.Text = [char]13
